# Adds three new "Attendance" log sheets after Sheet1, each a snapshot of the
# attendance log at a later point in time (rows accumulate / the last
# timestamp of the final entry gets corrected across snapshots).

$wb = $excel.ActiveWorkbook

# Reference the current header (Name/Date/Timestamp) formatting on Sheet1 so
# the new sheets' header rows pick up the same bold + bordered + centered
# style (style index 1) instead of inventing a new one.
$headerSrc = $wb.Worksheets.Item(1).Range("A1:C1")

function Add-AttendanceSheet($sheetName, $rows) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $sheetName

    # Copy header formatting (font/border/alignment) from Sheet1, formats only.
    $headerSrc.Copy()
    $ws.Range("A1:C1").PasteSpecial(-4122)

    $ws.Range("A1").Value = "Name"
    $ws.Range("B1").Value = "Date"
    $ws.Range("C1").Value = "Timestamp"

    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]

        # Dates/timestamps are plain log text, not real Excel dates/times -
        # force text entry (via a temporary "@" text format) so they don't
        # get auto-converted into date/time serial numbers, then drop the
        # temporary format so the cell is left with the default style
        # (no explicit "s" attribute), matching the rest of the data rows.
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 2).ClearFormats()

        $ws.Cells.Item($r, 3).NumberFormat = "@"
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 3).ClearFormats()

        $r = $r + 1
    }
}

$attendance = @(
    ,@("owen-4",    "2023-04-25", "23:12:35")
    ,@("owen-4",    "2023-04-25", "23:13:26")
    ,@("owen-4",    "2023-04-25", "23:14:08")
    ,@("Noom-157",  "2023-04-26", "01:04:13")
    ,@("Noom-157",  "2023-04-26", "01:10:14")
    ,@("owen-4",    "2023-04-26", "20:37:39")
    ,@("Noom-157",  "2023-04-26", "21:50:37")
)

$attendance1 = @(
    ,@("owen-4",    "2023-04-25", "23:12:35")
    ,@("owen-4",    "2023-04-25", "23:13:26")
    ,@("owen-4",    "2023-04-25", "23:14:08")
    ,@("Noom-157",  "2023-04-26", "01:04:13")
    ,@("Noom-157",  "2023-04-26", "01:10:14")
    ,@("owen-4",    "2023-04-26", "20:37:39")
    ,@("Noom-157",  "2023-04-26", "21:56:26")
)

$attendance2 = @(
    ,@("owen-4",          "2023-04-25", "23:12:35")
    ,@("owen-4",          "2023-04-25", "23:13:26")
    ,@("owen-4",          "2023-04-25", "23:14:08")
    ,@("Noom-157",        "2023-04-26", "01:04:13")
    ,@("Noom-157",        "2023-04-26", "01:10:14")
    ,@("owen-4",          "2023-04-26", "20:37:39")
    ,@("Noom-157",        "2023-04-26", "21:56:26")
    ,@("owen-64070103",   "2023-04-26", "21:56:30")
)

Add-AttendanceSheet "Attendance" $attendance
Add-AttendanceSheet "Attendance1" $attendance1
Add-AttendanceSheet "Attendance2" $attendance2
